$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap shuffled match rows back based on canonical order (same-date groups) ---
# Row 3
$ws.Cells.Item(3, 6).Value = "Spisska Nova Ves"
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = "Presov"
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 4.59
$ws.Cells.Item(3, 11).Value = "28/07/2023 03:42"
$ws.Cells.Item(3, 12).Value = 5.61
$ws.Cells.Item(3, 13).Value = "29/07/2023 15:53"
$ws.Cells.Item(3, 14).Value = 3.99
$ws.Cells.Item(3, 15).Value = "28/07/2023 03:42"
$ws.Cells.Item(3, 16).Value = 4.81
$ws.Cells.Item(3, 17).Value = "29/07/2023 15:53"
$ws.Cells.Item(3, 18).Value = 1.57
$ws.Cells.Item(3, 19).Value = "28/07/2023 03:42"
$ws.Cells.Item(3, 20).Value = 1.46
$ws.Cells.Item(3, 21).Value = "29/07/2023 15:53"
$ws.Cells.Item(3, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-presov/tnW4iUs2/"

# Row 7
$ws.Cells.Item(7, 6).Value = "Petrzalka"
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = "Pohronie"
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 2.3
$ws.Cells.Item(7, 11).Value = "28/07/2023 05:13"
$ws.Cells.Item(7, 12).Value = 2.42
$ws.Cells.Item(7, 13).Value = "28/07/2023 20:04"
$ws.Cells.Item(7, 14).Value = 3.3
$ws.Cells.Item(7, 15).Value = "28/07/2023 05:13"
$ws.Cells.Item(7, 16).Value = 3.37
$ws.Cells.Item(7, 17).Value = "29/07/2023 15:03"
$ws.Cells.Item(7, 18).Value = 2.73
$ws.Cells.Item(7, 19).Value = "28/07/2023 05:13"
$ws.Cells.Item(7, 20).Value = 2.72
$ws.Cells.Item(7, 21).Value = "29/07/2023 09:34"
$ws.Cells.Item(7, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-fk-pohronie/O8UXf3lf/"

# Row 75
$ws.Cells.Item(75, 6).Value = "Puchov"
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = "L. Mikulas"
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 2.16
$ws.Cells.Item(75, 11).Value = "29/09/2023 02:42"
$ws.Cells.Item(75, 12).Value = 2.4
$ws.Cells.Item(75, 13).Value = "30/09/2023 15:28"
$ws.Cells.Item(75, 14).Value = 3.39
$ws.Cells.Item(75, 15).Value = "29/09/2023 02:42"
$ws.Cells.Item(75, 16).Value = 3.45
$ws.Cells.Item(75, 17).Value = "30/09/2023 15:28"
$ws.Cells.Item(75, 18).Value = 2.81
$ws.Cells.Item(75, 19).Value = "29/09/2023 02:42"
$ws.Cells.Item(75, 20).Value = 2.74
$ws.Cells.Item(75, 21).Value = "30/09/2023 15:28"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-l-mikulas/WlQ0eh1Q/"

# Row 76
$ws.Cells.Item(76, 6).Value = "D. Kubin"
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = "Malzenice"
$ws.Cells.Item(76, 9).Value = 3
$ws.Cells.Item(76, 10).Value = 2.02
$ws.Cells.Item(76, 11).Value = "29/09/2023 02:42"
$ws.Cells.Item(76, 12).Value = 2
$ws.Cells.Item(76, 13).Value = "30/09/2023 15:22"
$ws.Cells.Item(76, 14).Value = 3.36
$ws.Cells.Item(76, 15).Value = "29/09/2023 02:42"
$ws.Cells.Item(76, 16).Value = 3.61
$ws.Cells.Item(76, 17).Value = "30/09/2023 15:22"
$ws.Cells.Item(76, 18).Value = 3.1
$ws.Cells.Item(76, 19).Value = "29/09/2023 02:42"
$ws.Cells.Item(76, 20).Value = 3.39
$ws.Cells.Item(76, 21).Value = "30/09/2023 15:22"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-malzenice/tOFzkjvf/"

# Row 77
$ws.Cells.Item(77, 6).Value = "Slovan Bratislava B"
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = "Pohronie"
$ws.Cells.Item(77, 9).Value = 3
$ws.Cells.Item(77, 10).Value = 2.23
$ws.Cells.Item(77, 11).Value = "29/09/2023 21:42"
$ws.Cells.Item(77, 12).Value = 2.37
$ws.Cells.Item(77, 13).Value = "01/10/2023 10:24"
$ws.Cells.Item(77, 14).Value = 3.32
$ws.Cells.Item(77, 15).Value = "29/09/2023 21:42"
$ws.Cells.Item(77, 16).Value = 3.63
$ws.Cells.Item(77, 17).Value = "01/10/2023 10:24"
$ws.Cells.Item(77, 18).Value = 2.74
$ws.Cells.Item(77, 19).Value = "29/09/2023 21:42"
$ws.Cells.Item(77, 20).Value = 2.67
$ws.Cells.Item(77, 21).Value = "01/10/2023 10:24"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-fk-pohronie/2w4qmU86/"

# Row 78
$ws.Cells.Item(78, 6).Value = "Samorin"
$ws.Cells.Item(78, 7).Value = 1
$ws.Cells.Item(78, 8).Value = "Komarno"
$ws.Cells.Item(78, 9).Value = 4
$ws.Cells.Item(78, 10).Value = 3.75
$ws.Cells.Item(78, 11).Value = "29/09/2023 21:42"
$ws.Cells.Item(78, 12).Value = 4.36
$ws.Cells.Item(78, 13).Value = "01/10/2023 10:22"
$ws.Cells.Item(78, 14).Value = 3.65
$ws.Cells.Item(78, 15).Value = "29/09/2023 21:42"
$ws.Cells.Item(78, 16).Value = 3.88
$ws.Cells.Item(78, 17).Value = "01/10/2023 10:22"
$ws.Cells.Item(78, 18).Value = 1.74
$ws.Cells.Item(78, 19).Value = "29/09/2023 21:42"
$ws.Cells.Item(78, 20).Value = 1.71
$ws.Cells.Item(78, 21).Value = "01/10/2023 10:22"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/samorin-komarno/On5ulAg0/"

# Row 79
$ws.Cells.Item(79, 6).Value = "Petrzalka"
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = "FK Humenne"
$ws.Cells.Item(79, 9).Value = 1
$ws.Cells.Item(79, 10).Value = 1.73
$ws.Cells.Item(79, 11).Value = "29/09/2023 21:42"
$ws.Cells.Item(79, 12).Value = 1.65
$ws.Cells.Item(79, 13).Value = "01/10/2023 10:26"
$ws.Cells.Item(79, 14).Value = 3.63
$ws.Cells.Item(79, 15).Value = "29/09/2023 21:42"
$ws.Cells.Item(79, 16).Value = 4.1
$ws.Cells.Item(79, 17).Value = "01/10/2023 10:28"
$ws.Cells.Item(79, 18).Value = 3.82
$ws.Cells.Item(79, 19).Value = "29/09/2023 21:42"
$ws.Cells.Item(79, 20).Value = 4.53
$ws.Cells.Item(79, 21).Value = "01/10/2023 10:26"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-fk-humenne/dEInnlOC/"

# Row 91
$ws.Cells.Item(91, 6).Value = "FK Humenne"
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = "Malzenice"
$ws.Cells.Item(91, 9).Value = 1
$ws.Cells.Item(91, 10).Value = 1.53
$ws.Cells.Item(91, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(91, 12).Value = 1.39
$ws.Cells.Item(91, 13).Value = "14/10/2023 14:39"
$ws.Cells.Item(91, 14).Value = 3.95
$ws.Cells.Item(91, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(91, 16).Value = 4.47
$ws.Cells.Item(91, 17).Value = "14/10/2023 14:39"
$ws.Cells.Item(91, 18).Value = 4.82
$ws.Cells.Item(91, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(91, 20).Value = 7.94
$ws.Cells.Item(91, 21).Value = "14/10/2023 14:39"
$ws.Cells.Item(91, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-malzenice/z3np52Ui/"

# Row 92
$ws.Cells.Item(92, 6).Value = "Povazska Bystrica"
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 8).Value = "Zilina B"
$ws.Cells.Item(92, 9).Value = 1
$ws.Cells.Item(92, 10).Value = 1.77
$ws.Cells.Item(92, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(92, 12).Value = 2.15
$ws.Cells.Item(92, 13).Value = "14/10/2023 14:56"
$ws.Cells.Item(92, 14).Value = 3.81
$ws.Cells.Item(92, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(92, 16).Value = 3.74
$ws.Cells.Item(92, 17).Value = "14/10/2023 14:56"
$ws.Cells.Item(92, 18).Value = 3.47
$ws.Cells.Item(92, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(92, 20).Value = 2.95
$ws.Cells.Item(92, 21).Value = "14/10/2023 14:56"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-zilina/dCmt6rFo/"

# Row 93
$ws.Cells.Item(93, 6).Value = "D. Kubin"
$ws.Cells.Item(93, 7).Value = 2
$ws.Cells.Item(93, 8).Value = "Spisska Nova Ves"
$ws.Cells.Item(93, 9).Value = 1
$ws.Cells.Item(93, 10).Value = 2
$ws.Cells.Item(93, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(93, 12).Value = 2.45
$ws.Cells.Item(93, 13).Value = "14/10/2023 14:48"
$ws.Cells.Item(93, 14).Value = 3.36
$ws.Cells.Item(93, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(93, 16).Value = 3.49
$ws.Cells.Item(93, 17).Value = "14/10/2023 14:51"
$ws.Cells.Item(93, 18).Value = 3.15
$ws.Cells.Item(93, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(93, 20).Value = 2.65
$ws.Cells.Item(93, 21).Value = "14/10/2023 14:48"
$ws.Cells.Item(93, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-spisska-nova-ves/IRyk4Mqb/"

# Row 94
$ws.Cells.Item(94, 6).Value = "Puchov"
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = "Komarno"
$ws.Cells.Item(94, 9).Value = 1
$ws.Cells.Item(94, 10).Value = 2.69
$ws.Cells.Item(94, 11).Value = "13/10/2023 02:13"
$ws.Cells.Item(94, 12).Value = 2.36
$ws.Cells.Item(94, 13).Value = "14/10/2023 14:52"
$ws.Cells.Item(94, 14).Value = 3.23
$ws.Cells.Item(94, 15).Value = "13/10/2023 02:13"
$ws.Cells.Item(94, 16).Value = 3.38
$ws.Cells.Item(94, 17).Value = "14/10/2023 14:52"
$ws.Cells.Item(94, 18).Value = 2.31
$ws.Cells.Item(94, 19).Value = "13/10/2023 02:13"
$ws.Cells.Item(94, 20).Value = 2.84
$ws.Cells.Item(94, 21).Value = "14/10/2023 14:52"
$ws.Cells.Item(94, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-komarno/8dfbsaxo/"

# Row 100
$ws.Cells.Item(100, 6).Value = "Presov"
$ws.Cells.Item(100, 7).Value = 1
$ws.Cells.Item(100, 8).Value = "Petrzalka"
$ws.Cells.Item(100, 9).Value = 1
$ws.Cells.Item(100, 10).Value = 1.56
$ws.Cells.Item(100, 11).Value = "20/10/2023 01:42"
$ws.Cells.Item(100, 12).Value = 1.85
$ws.Cells.Item(100, 13).Value = "21/10/2023 14:19"
$ws.Cells.Item(100, 14).Value = 3.84
$ws.Cells.Item(100, 15).Value = "20/10/2023 01:42"
$ws.Cells.Item(100, 16).Value = 3.74
$ws.Cells.Item(100, 17).Value = "21/10/2023 14:19"
$ws.Cells.Item(100, 18).Value = 4.69
$ws.Cells.Item(100, 19).Value = "20/10/2023 01:42"
$ws.Cells.Item(100, 20).Value = 3.82
$ws.Cells.Item(100, 21).Value = "21/10/2023 14:19"
$ws.Cells.Item(100, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/presov-petrzalka/Opf2abbT/"

# Row 101
$ws.Cells.Item(101, 6).Value = "Spisska Nova Ves"
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(101, 8).Value = "FK Humenne"
$ws.Cells.Item(101, 9).Value = 3
$ws.Cells.Item(101, 10).Value = 3.02
$ws.Cells.Item(101, 11).Value = "20/10/2023 01:42"
$ws.Cells.Item(101, 12).Value = 3.43
$ws.Cells.Item(101, 13).Value = "21/10/2023 14:27"
$ws.Cells.Item(101, 14).Value = 3.19
$ws.Cells.Item(101, 15).Value = "20/10/2023 01:42"
$ws.Cells.Item(101, 16).Value = 3.47
$ws.Cells.Item(101, 17).Value = "21/10/2023 14:27"
$ws.Cells.Item(101, 18).Value = 2.13
$ws.Cells.Item(101, 19).Value = "20/10/2023 01:42"
$ws.Cells.Item(101, 20).Value = 2.04
$ws.Cells.Item(101, 21).Value = "21/10/2023 14:27"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-fk-humenne/fTpYgxMj/"

# Row 102
$ws.Cells.Item(102, 6).Value = "Malzenice"
$ws.Cells.Item(102, 7).Value = 2
$ws.Cells.Item(102, 8).Value = "Myjava"
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 10).Value = 2.59
$ws.Cells.Item(102, 11).Value = "20/10/2023 01:42"
$ws.Cells.Item(102, 12).Value = 3.04
$ws.Cells.Item(102, 13).Value = "21/10/2023 14:28"
$ws.Cells.Item(102, 14).Value = 3.23
$ws.Cells.Item(102, 15).Value = "20/10/2023 01:42"
$ws.Cells.Item(102, 16).Value = 3.14
$ws.Cells.Item(102, 17).Value = "21/10/2023 14:28"
$ws.Cells.Item(102, 18).Value = 2.4
$ws.Cells.Item(102, 19).Value = "20/10/2023 01:42"
$ws.Cells.Item(102, 20).Value = 2.36
$ws.Cells.Item(102, 21).Value = "21/10/2023 14:23"
$ws.Cells.Item(102, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-myjava/4bSOGaT3/"

# Row 105
$ws.Cells.Item(105, 6).Value = "FK Humenne"
$ws.Cells.Item(105, 7).Value = 4
$ws.Cells.Item(105, 8).Value = "L. Mikulas"
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(105, 10).Value = 2.3
$ws.Cells.Item(105, 11).Value = "27/10/2023 02:42"
$ws.Cells.Item(105, 12).Value = 1.96
$ws.Cells.Item(105, 13).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 14).Value = 3.18
$ws.Cells.Item(105, 15).Value = "27/10/2023 02:42"
$ws.Cells.Item(105, 16).Value = 3.59
$ws.Cells.Item(105, 17).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 18).Value = 2.74
$ws.Cells.Item(105, 19).Value = "27/10/2023 02:42"
$ws.Cells.Item(105, 20).Value = 3.52
$ws.Cells.Item(105, 21).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-l-mikulas/x6FB7ejj/"

# Row 106
$ws.Cells.Item(106, 6).Value = "Povazska Bystrica"
$ws.Cells.Item(106, 7).Value = 3
$ws.Cells.Item(106, 8).Value = "Slovan Bratislava B"
$ws.Cells.Item(106, 9).Value = 1
$ws.Cells.Item(106, 10).Value = 1.8
$ws.Cells.Item(106, 11).Value = "27/10/2023 02:42"
$ws.Cells.Item(106, 12).Value = 1.46
$ws.Cells.Item(106, 13).Value = "28/10/2023 13:57"
$ws.Cells.Item(106, 14).Value = 3.56
$ws.Cells.Item(106, 15).Value = "27/10/2023 02:42"
$ws.Cells.Item(106, 16).Value = 4.44
$ws.Cells.Item(106, 17).Value = "28/10/2023 13:57"
$ws.Cells.Item(106, 18).Value = 3.58
$ws.Cells.Item(106, 19).Value = "27/10/2023 02:42"
$ws.Cells.Item(106, 20).Value = 6.15
$ws.Cells.Item(106, 21).Value = "28/10/2023 13:57"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-slovan-bratislava/ObZzEcDM/"

# Row 107
$ws.Cells.Item(107, 6).Value = "Puchov"
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = "Pohronie"
$ws.Cells.Item(107, 9).Value = 4
$ws.Cells.Item(107, 10).Value = 1.75
$ws.Cells.Item(107, 11).Value = "27/10/2023 02:42"
$ws.Cells.Item(107, 12).Value = 1.71
$ws.Cells.Item(107, 13).Value = "28/10/2023 14:21"
$ws.Cells.Item(107, 14).Value = 3.6
$ws.Cells.Item(107, 15).Value = "27/10/2023 02:42"
$ws.Cells.Item(107, 16).Value = 3.81
$ws.Cells.Item(107, 17).Value = "28/10/2023 14:21"
$ws.Cells.Item(107, 18).Value = 3.77
$ws.Cells.Item(107, 19).Value = "27/10/2023 02:42"
$ws.Cells.Item(107, 20).Value = 4.45
$ws.Cells.Item(107, 21).Value = "28/10/2023 14:21"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-fk-pohronie/QDGXEwbG/"

# Row 108
$ws.Cells.Item(108, 6).Value = "D. Kubin"
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = "Presov"
$ws.Cells.Item(108, 9).Value = 3
$ws.Cells.Item(108, 10).Value = 5.32
$ws.Cells.Item(108, 11).Value = "27/10/2023 02:42"
$ws.Cells.Item(108, 12).Value = 8.949999999999999
$ws.Cells.Item(108, 13).Value = "28/10/2023 14:28"
$ws.Cells.Item(108, 14).Value = 3.95
$ws.Cells.Item(108, 15).Value = "27/10/2023 02:42"
$ws.Cells.Item(108, 16).Value = 5.83
$ws.Cells.Item(108, 17).Value = "28/10/2023 14:28"
$ws.Cells.Item(108, 18).Value = 1.48
$ws.Cells.Item(108, 19).Value = "27/10/2023 02:42"
$ws.Cells.Item(108, 20).Value = 1.26
$ws.Cells.Item(108, 21).Value = "28/10/2023 14:28"
$ws.Cells.Item(108, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-presov/ILYvDHSS/"

# Row 109
$ws.Cells.Item(109, 6).Value = "Petrzalka"
$ws.Cells.Item(109, 7).Value = 3
$ws.Cells.Item(109, 8).Value = "Komarno"
$ws.Cells.Item(109, 9).Value = 1
$ws.Cells.Item(109, 10).Value = 2.89
$ws.Cells.Item(109, 11).Value = "27/10/2023 23:42"
$ws.Cells.Item(109, 12).Value = 2.91
$ws.Cells.Item(109, 13).Value = "29/10/2023 10:02"
$ws.Cells.Item(109, 14).Value = 3.2
$ws.Cells.Item(109, 15).Value = "27/10/2023 23:42"
$ws.Cells.Item(109, 16).Value = 3.44
$ws.Cells.Item(109, 17).Value = "29/10/2023 10:02"
$ws.Cells.Item(109, 18).Value = 2.2
$ws.Cells.Item(109, 19).Value = "27/10/2023 23:42"
$ws.Cells.Item(109, 20).Value = 2.29
$ws.Cells.Item(109, 21).Value = "29/10/2023 10:02"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-komarno/GdEF6F5d/"

# Row 110
$ws.Cells.Item(110, 6).Value = "Malzenice"
$ws.Cells.Item(110, 7).Value = 3
$ws.Cells.Item(110, 8).Value = "Spisska Nova Ves"
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = 1.94
$ws.Cells.Item(110, 11).Value = "28/10/2023 08:13"
$ws.Cells.Item(110, 12).Value = 2.03
$ws.Cells.Item(110, 13).Value = "29/10/2023 10:21"
$ws.Cells.Item(110, 14).Value = 3.34
$ws.Cells.Item(110, 15).Value = "28/10/2023 08:13"
$ws.Cells.Item(110, 16).Value = 3.33
$ws.Cells.Item(110, 17).Value = "29/10/2023 10:29"
$ws.Cells.Item(110, 18).Value = 3.32
$ws.Cells.Item(110, 19).Value = "28/10/2023 08:13"
$ws.Cells.Item(110, 20).Value = 3.59
$ws.Cells.Item(110, 21).Value = "29/10/2023 10:21"
$ws.Cells.Item(110, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-spisska-nova-ves/KUK68yyp/"

# Row 111
$ws.Cells.Item(111, 6).Value = "Samorin"
$ws.Cells.Item(111, 7).Value = 2
$ws.Cells.Item(111, 8).Value = "Myjava"
$ws.Cells.Item(111, 9).Value = 2
$ws.Cells.Item(111, 10).Value = 2.32
$ws.Cells.Item(111, 11).Value = "27/10/2023 23:42"
$ws.Cells.Item(111, 12).Value = 2.78
$ws.Cells.Item(111, 13).Value = "29/10/2023 09:48"
$ws.Cells.Item(111, 14).Value = 3.35
$ws.Cells.Item(111, 15).Value = "27/10/2023 23:42"
$ws.Cells.Item(111, 16).Value = 3.53
$ws.Cells.Item(111, 17).Value = "29/10/2023 09:48"
$ws.Cells.Item(111, 18).Value = 2.66
$ws.Cells.Item(111, 19).Value = "27/10/2023 23:42"
$ws.Cells.Item(111, 20).Value = 2.33
$ws.Cells.Item(111, 21).Value = "29/10/2023 09:48"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/samorin-myjava/Ua2cVbrc/"

# --- Append new rows 121-124 ---
# New row 121
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(121, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(121, 5).PasteSpecial(-4122)
$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = "slovakia"
$ws.Cells.Item(121, 3).Value = "2-liga"
$ws.Cells.Item(121, 4).Value = "2023-2024"
$ws.Cells.Item(121, 5).Value = 45241.4375
$ws.Cells.Item(121, 6).Value = "L. Mikulas"
$ws.Cells.Item(121, 7).Value = 4
$ws.Cells.Item(121, 8).Value = "Malzenice"
$ws.Cells.Item(121, 9).Value = 2
$ws.Cells.Item(121, 10).Value = 1.48
$ws.Cells.Item(121, 11).Value = "09/11/2023 22:42"
$ws.Cells.Item(121, 12).Value = 1.49
$ws.Cells.Item(121, 13).Value = "11/11/2023 10:15"
$ws.Cells.Item(121, 14).Value = 4.08
$ws.Cells.Item(121, 15).Value = "09/11/2023 22:42"
$ws.Cells.Item(121, 16).Value = 4.58
$ws.Cells.Item(121, 17).Value = "11/11/2023 10:27"
$ws.Cells.Item(121, 18).Value = 5.08
$ws.Cells.Item(121, 19).Value = "09/11/2023 22:42"
$ws.Cells.Item(121, 20).Value = 5.56
$ws.Cells.Item(121, 21).Value = "11/11/2023 10:27"
$ws.Cells.Item(121, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/l-mikulas-malzenice/h2FmtGdF/"

# New row 122
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(122, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(122, 5).PasteSpecial(-4122)
$ws.Cells.Item(122, 1).Value = 121
$ws.Cells.Item(122, 2).Value = "slovakia"
$ws.Cells.Item(122, 3).Value = "2-liga"
$ws.Cells.Item(122, 4).Value = "2023-2024"
$ws.Cells.Item(122, 5).Value = 45241.4375
$ws.Cells.Item(122, 6).Value = "Slovan Bratislava B"
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = "Trebisov"
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 1.74
$ws.Cells.Item(122, 11).Value = "09/11/2023 22:42"
$ws.Cells.Item(122, 12).Value = 2.09
$ws.Cells.Item(122, 13).Value = "11/11/2023 10:10"
$ws.Cells.Item(122, 14).Value = 3.66
$ws.Cells.Item(122, 15).Value = "09/11/2023 22:42"
$ws.Cells.Item(122, 16).Value = 3.64
$ws.Cells.Item(122, 17).Value = "11/11/2023 10:17"
$ws.Cells.Item(122, 18).Value = 3.87
$ws.Cells.Item(122, 19).Value = "09/11/2023 22:42"
$ws.Cells.Item(122, 20).Value = 3.1
$ws.Cells.Item(122, 21).Value = "11/11/2023 10:10"
$ws.Cells.Item(122, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-trebisov/lAfOZDZr/"

# New row 123
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(123, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(123, 5).PasteSpecial(-4122)
$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = "slovakia"
$ws.Cells.Item(123, 3).Value = "2-liga"
$ws.Cells.Item(123, 4).Value = "2023-2024"
$ws.Cells.Item(123, 5).Value = 45241.54166666666
$ws.Cells.Item(123, 6).Value = "Presov"
$ws.Cells.Item(123, 7).Value = 2
$ws.Cells.Item(123, 8).Value = "Spisska Nova Ves"
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 1.18
$ws.Cells.Item(123, 11).Value = "10/11/2023 01:13"
$ws.Cells.Item(123, 12).Value = 1.17
$ws.Cells.Item(123, 13).Value = "11/11/2023 12:57"
$ws.Cells.Item(123, 14).Value = 6.08
$ws.Cells.Item(123, 15).Value = "10/11/2023 01:13"
$ws.Cells.Item(123, 16).Value = 6.84
$ws.Cells.Item(123, 17).Value = "11/11/2023 12:57"
$ws.Cells.Item(123, 18).Value = 9.97
$ws.Cells.Item(123, 19).Value = "10/11/2023 01:13"
$ws.Cells.Item(123, 20).Value = 16.39
$ws.Cells.Item(123, 21).Value = "11/11/2023 12:57"
$ws.Cells.Item(123, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/presov-spisska-nova-ves/4xJiuzCL/"

# New row 124
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(124, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(124, 5).PasteSpecial(-4122)
$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = "slovakia"
$ws.Cells.Item(124, 3).Value = "2-liga"
$ws.Cells.Item(124, 4).Value = "2023-2024"
$ws.Cells.Item(124, 5).Value = 45241.54166666666
$ws.Cells.Item(124, 6).Value = "Komarno"
$ws.Cells.Item(124, 7).Value = 1
$ws.Cells.Item(124, 8).Value = "FK Humenne"
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = 1.56
$ws.Cells.Item(124, 11).Value = "10/11/2023 01:13"
$ws.Cells.Item(124, 12).Value = 1.67
$ws.Cells.Item(124, 13).Value = "11/11/2023 12:45"
$ws.Cells.Item(124, 14).Value = 3.77
$ws.Cells.Item(124, 15).Value = "10/11/2023 01:13"
$ws.Cells.Item(124, 16).Value = 3.74
$ws.Cells.Item(124, 17).Value = "11/11/2023 12:45"
$ws.Cells.Item(124, 18).Value = 4.81
$ws.Cells.Item(124, 19).Value = "10/11/2023 01:13"
$ws.Cells.Item(124, 20).Value = 4.9
$ws.Cells.Item(124, 21).Value = "11/11/2023 12:45"
$ws.Cells.Item(124, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/komarno-fk-humenne/buGqsds9/"

$ws.Application.CutCopyMode = $false